$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.542.57'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.41'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.35'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5144'
$ws.Range("E7").Value = '  -3.62%  '
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08451'
$ws.Range("E9").Value = '  +8.11%  '
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.398'
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.02'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.501'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.812.34'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001136'
$ws.Range("E17").Value = '  +3.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.80'
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06688'
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.75'
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.086'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.559.81'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.43'
$ws.Range("E24").Value = '  +2.12%  '
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.10'
$ws.Range("E26").Value = '  +1.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.86'
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.023.48'
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.411'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.04'
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E31").Value = '  -4.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1080'
$ws.Range("E32").Value = '  -3.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.750'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07599'
$ws.Range("E34").Value = '  +4.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.683'
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2225'
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02367'
$ws.Range("E37").Value = '  +0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.202'
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.739'
$ws.Range("E39").Value = '  -2.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6331'
$ws.Range("E40").Value = '  +0.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.25'
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.401'
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.765'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5912'
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '126.00'
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.200'
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06983'
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("E51").Value = '  -0.72%  '
